$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title texts (Volume/Number and date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "40"

$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "9/29/2025"
$c9.Characters(47, 9).Text = "10/5/2025"

# --- Update crime statistics table (rows 14-30) ---

# Row 14
$ws.Range("I14").Value = 8
$ws.Range("K14").Value = 60
$ws.Range("L14").Value = -20
$ws.Range("M14").Value = -11.111111111111
$ws.Range("N14").Value = -68

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("L15").Value = 20
$ws.Range("M15").Value = 130.769230769231
$ws.Range("N15").Value = -51.612903225806

# Row 16
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = 19.354838709677
$ws.Range("I16").Value = 349
$ws.Range("J16").Value = 354
$ws.Range("K16").Value = -1.412429378531
$ws.Range("L16").Value = 9.404388714733
$ws.Range("M16").Value = 55.803571428571
$ws.Range("N16").Value = -63.185654008438

# Row 17
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 15.384615384615
$ws.Range("F17").Value = 66
$ws.Range("G17").Value = 57
$ws.Range("H17").Value = 15.789473684210
$ws.Range("I17").Value = 607
$ws.Range("J17").Value = 579
$ws.Range("K17").Value = 4.835924006908
$ws.Range("L17").Value = 10.163339382940
$ws.Range("M17").Value = 152.916666666667
$ws.Range("N17").Value = -20.341207349081

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -28
$ws.Range("I18").Value = 181
$ws.Range("J18").Value = 221
$ws.Range("K18").Value = -18.099547511312
$ws.Range("L18").Value = 20.666666666666
$ws.Range("M18").Value = 54.700854700854
$ws.Range("N18").Value = -78.400954653937

# Row 19
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 40
$ws.Range("I19").Value = 460
$ws.Range("J19").Value = 423
$ws.Range("K19").Value = 8.747044917257
$ws.Range("L19").Value = 16.455696202531
$ws.Range("M19").Value = 148.648648648649
$ws.Range("N19").Value = 64.874551971326

# Row 20
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 4.545454545454
$ws.Range("I20").Value = 199
$ws.Range("J20").Value = 212
$ws.Range("K20").Value = -6.132075471698
$ws.Range("L20").Value = -42.651296829971
$ws.Range("M20").Value = 128.735632183908
$ws.Range("N20").Value = -45.628415300546

# Row 21
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 15.789473684210
$ws.Range("F21").Value = 196
$ws.Range("G21").Value = 170
$ws.Range("H21").Value = 15.294117647058
$ws.Range("I21").Value = 1834
$ws.Range("J21").Value = 1826
$ws.Range("K21").Value = 0.438116100766
$ws.Range("L21").Value = 2.058987200890
$ws.Range("M21").Value = 109.6
$ws.Range("N21").Value = -44.085365853658

# Row 23
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 27
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = 17.391304347826
$ws.Range("I23").Value = 234
$ws.Range("J23").Value = 329
$ws.Range("K23").Value = -28.875379939209
$ws.Range("L23").Value = -27.554179566563
$ws.Range("M23").Value = 58.108108108108

# Row 24
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 20.833333333333
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = -12.790697674418
$ws.Range("I24").Value = 909
$ws.Range("J24").Value = 781
$ws.Range("K24").Value = 16.389244558258
$ws.Range("L24").Value = 0.441988950276
$ws.Range("M24").Value = 54.591836734693

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 167
$ws.Range("J25").Value = 133
$ws.Range("K25").Value = 25.563909774436
$ws.Range("L25").Value = -18.137254901960

# Row 26
$ws.Range("C26").Value = 22
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = 4.761904761904
$ws.Range("F26").Value = 76
$ws.Range("G26").Value = 87
$ws.Range("H26").Value = -12.643678160919
$ws.Range("I26").Value = 705
$ws.Range("J26").Value = 753
$ws.Range("K26").Value = -6.374501992031
$ws.Range("L26").Value = -14.233576642335
$ws.Range("M26").Value = 4.754829123328

# Row 27
$ws.Range("F27").Value = 3
$ws.Range("L27").Value = -20.930232558139

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 44
$ws.Range("J28").Value = 65
$ws.Range("K28").Value = -32.307692307692
$ws.Range("L28").Value = -41.333333333333

# Row 29
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 29
$ws.Range("J29").Value = 38
$ws.Range("K29").Value = -23.684210526315
$ws.Range("L29").Value = -14.705882352941
$ws.Range("M29").Value = -25.641025641025
$ws.Range("N29").Value = -65.476190476190

# Row 30
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 27
$ws.Range("J30").Value = 33
$ws.Range("K30").Value = -18.181818181818
$ws.Range("L30").Value = -3.571428571428
$ws.Range("M30").Value = -18.181818181818
$ws.Range("N30").Value = -67.857142857142

# --- Cells that switch between numeric and text representation ---
# (value is set first, then the number format is copied from a donor
#  cell that already carries the correct target style, so the cells
#  style index matches the target workbook exactly)

# C14: s:'0' -> n:'1'  (style -> 14, donor I14)
$ws.Range("C14").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null

# F14: s:'0' -> n:'1'  (style -> 14, donor I14)
$ws.Range("F14").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null

# G15: n:'3' -> s:'0'  (style -> 13, donor D14)
$ws.Range("G15").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null

# H15: n:'100' -> s:'***.*'  (style -> 13, donor D14)
$ws.Range("H15").Value = "'***.*"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null

# C27: n:'1' -> s:'0'  (style -> 13, donor D14)
$ws.Range("C27").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null

# G27: n:'3' -> s:'0'  (style -> 13, donor D14)
$ws.Range("G27").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null

# H27: n:'133.333333333333' -> s:'***.*'  (style -> 13, donor D14)
$ws.Range("H27").Value = "'***.*"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null

# C29: s:'0' -> n:'1'  (style -> 14, donor I14)
$ws.Range("C29").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null

# D29: s:'0' -> n:'2'  (style -> 14, donor I14)
$ws.Range("D29").Value = 2
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null

# E29: s:'***.*' -> n:'-50'  (style -> 15, donor K14)
$ws.Range("E29").Value = -50
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null

# C30: s:'0' -> n:'1'  (style -> 14, donor I14)
$ws.Range("C30").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null

# D30: s:'0' -> n:'2'  (style -> 14, donor I14)
$ws.Range("D30").Value = 2
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null

# E30: s:'***.*' -> n:'-50'  (style -> 15, donor K14)
$ws.Range("E30").Value = -50
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
